# Edit script: updates the "LojaVirtual-ProcessoDevolucao" worksheet
# with the additional test-documentation rows (TCID CT-ValProsTroTrl001
# "Loja Virtual" flow), a mailto hyperlink on G8, the widened column G,
# and moves the active selection back to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LojaVirtual-ProcessoDevolucao")

$ws.Range("F7").Value = "Preencher o campo Pedido"
$ws.Range("G7").Value = 318
$ws.Range("F8").Value = "Preencher o campo E-mail utilizado na compra"
$ws.Range("F9").Value = "Clicar no botão Buscar e continuar"
$ws.Range("H9").Value = "Deve prosseguir para a próxima tela."
$ws.Range("F10").Value = "Na tela carregada, clicar no checkbox Selecionar todos"
$ws.Range("F11").Value = "No primeiro produto:"
$ws.Range("F12").Value = "Selecionar a quantidade"
$ws.Range("G12").Value = 1
$ws.Range("F13").Value = "Selecionar a Ação"
$ws.Range("G13").Value = "Trocar"
$ws.Range("F14").Value = "Selecionar o Motivo"
$ws.Range("G14").Value = "Defeito na troca"
$ws.Range("F15").Value = "Preencher o campo Como podemos Resolver?"
$ws.Range("G15").Value = "Digite um texto"
$ws.Range("F16").Value = "No próximo produto"
$ws.Range("F17").Value = "Selecionar a quantidade"
$ws.Range("G17").Value = 1
$ws.Range("F18").Value = "Selecionar a Ação"
$ws.Range("G18").Value = "Devolver"
$ws.Range("F19").Value = "Selecionar o Motivo"
$ws.Range("G19").Value = "Me arrependi"
$ws.Range("F20").Value = "Preencher o campo Como podemos Resolver?"
$ws.Range("G20").Value = "Digite outro texto"
$ws.Range("F21").Value = "Clicar no botão Continuar"
$ws.Range("H21").Value = "Deve prosseguir para a próxima tela."
$ws.Range("F22").Value = "Na tela carregada, deve apresentar a oferta de devolução em Vale-Compras"
$ws.Range("F23").Value = "Clicar no botão Vou esperar meu dinheiro por30 dias"
$ws.Range("H23").Value = "Deve prosseguir para a próxima tela."
$ws.Range("F24").Value = "Na tela carregada, deve apresentar a oferta de Vale-compras + Super oferta com timer"
$ws.Range("F25").Value = "Clicar no botão Deve prosseguir para a próxima tela. 30 dias"
$ws.Range("H25").Value = "Deve prosseguir para a próxima tela."
$ws.Range("F26").Value = "NA tela carregada, deve apresentar o titulo Detalhes Bancários"
$ws.Range("F27").Value = "Selecionar a opção"
$ws.Range("G27").Value = "Não tenho conta"
$ws.Range("G27").WrapText = $true
$ws.Range("F28").Value = "Clicar no botão Continuar"
$ws.Range("H28").Value = "Deve prosseguir para a próxima tela."
$ws.Range("F29").Value = "NA tela carregada, deve apresentar o titulo Selecione o método de devolução"
$ws.Range("F30").Value = "Validar se o endereço está correto"
$ws.Range("F31").Value = "Clicar no botão Selecione da Agência de Correios"
$ws.Range("F32").Value = "Clicar no botão Continuar"
$ws.Range("H32").Value = "Deve prosseguir para a próxima tela."
$ws.Range("F33").Value = "Conferir os dados da seção de Dados do Cliente"
$ws.Range("F34").Value = "Conferir os dados da seção de Produtos"
$ws.Range("F35").Value = "Marcar o checkbox de Li e concordo..."
$ws.Range("F36").Value = "Clicar no botão Continuar"
$ws.Range("H36").Value = "Deve prosseguir para a próxima tela."
$ws.Range("F37").Value = "Na tela seguinte, com a mensagem Sua solicitação foi realizada com sucesso!"
$ws.Range("F38").Value = "Clicar em uma nota"
$ws.Range("F39").Value = "Preencher o campo Comentário"
$ws.Range("G39").Value = "Qualquer texto"
$ws.Range("F40").Value = "Clicar em enviar avaliação"

# Hyperlink: tester@send4.com.br (mailto) placed on G8
$ws.Hyperlinks.Add($ws.Range("G8"), "mailto:tester@send4.com.br", "", "", "tester@send4.com.br")

# Column G needs to be wide enough to show the e-mail / numbers that were added
$ws.Columns.Item(7).ColumnWidth = 17

# Restore the active selection to C2 (it had drifted to C9)
$ws.Range("C2").Select()
